$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'56.208.55"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +3.29%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.479.60"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +2.00%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.03%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'488.89"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +5.11%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'147.21"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +12.04%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.997"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -0.27%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.510"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +3.77%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'2.487.23"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +2.52%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  +9.61%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.0972"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +2.57%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.333"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +5.98%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "'  +1.76%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'2.918.62"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +2.34%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'56.186.13"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +3.22%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'21.15"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +7.59%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("E17").Value = "'  +3.48%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'2.487.28"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +2.44%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'4.53"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +8.18%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'10.09"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +6.18%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'319.12"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +3.14%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.998"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +0.01%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'  +8.35%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'58.54"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +4.10%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'  +7.21%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  -0.84%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  +4.39%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'2.584.63"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +1.96%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'7.64"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +7.63%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'0.0₃0791"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +11.51%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  -0.21%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'149.09"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +1.96%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'18.23"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +3.13%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  +5.38%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "'  +4.35%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  +8.69%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'3.75"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +6.39%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.862"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +8.35%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  +3.63%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'3.52"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +8.86%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'  -0.24%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.0555"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +6.58%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.604"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +1.95%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'  +8.23%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = "'  +15.19%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("B46").Value = "'Bittensor"
$ws.Range("B46").Style = "Normal"
$ws.Range("C46").Value = "'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("C46").Style = "Normal"
$ws.Range("D46").Value = "'260.41"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +13.91%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("B47").Value = "'Stellar"
$ws.Range("B47").Style = "Normal"
$ws.Range("C47").Value = "'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("C47").Style = "Normal"
$ws.Range("D47").Value = "'0.0924"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +5.45%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("B48").Value = "'VeChain"
$ws.Range("B48").Style = "Normal"
$ws.Range("C48").Value = "'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("C48").Style = "Normal"
$ws.Range("D48").Value = "'0.0228"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +5.08%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("B49").Value = "'WhiteBITCoin"
$ws.Range("B49").Style = "Normal"
$ws.Range("C49").Value = "'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("C49").Style = "Normal"
$ws.Range("D49").Value = "'10.18"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +1.13%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'17.58"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +6.85%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'1.883.49"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -2.28%  "
$ws.Range("E51").Style = "Normal"
